$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

# Header date line
Replace-Text "2023-11-12 Sunday" "2023-11-13 Monday"

# Table row 1
Replace-Text "45÷6=7, 3" "66÷2=33, 0"
Replace-Text "76÷5=15, 1" "83÷9=9, 2"
Replace-Text "90÷9=10, 0" "71÷4=17, 3"
Replace-Text "34÷7=4, 6" "32÷6=5, 2"
Replace-Text "65÷6=10, 5" "33÷7=4, 5"

# Table row 2 (data row index 5, 1-based table row)
Replace-Text "46÷9=5, 1" "23÷2=11, 1"
Replace-Text "90÷2=45, 0" "47÷6=7, 5"
Replace-Text "72÷5=14, 2" "10÷7=1, 3"
Replace-Text "67÷3=22, 1" "38÷8=4, 6"
Replace-Text "47÷5=9, 2" "81÷6=13, 3"

# Table row 3 (data row index 9, 1-based table row)
Replace-Text "28÷5=5, 3" "91÷2=45, 1"
Replace-Text "95÷7=13, 4" "85÷6=14, 1"
Replace-Text "28÷2=14, 0" "80÷8=10, 0"
Replace-Text "85÷2=42, 1" "45÷7=6, 3"
Replace-Text "40÷5=8, 0" "48÷9=5, 3"

# Table row 4 (data row index 13, 1-based table row) - cells rearranged
$t = $d.Tables.Item(1)
$t.Cell(13,1).Range.Text = "50÷8=6, 2"
$t.Cell(13,2).Range.Text = "49÷3=16, 1"
$t.Cell(13,3).Range.Text = "21÷4=5, 1"
$t.Cell(13,4).Range.Text = "31÷8=3, 7"
$t.Cell(13,5).Range.Text = "37÷5=7, 2"

# Table row 5 (data row index 17, 1-based table row)
Replace-Text "32÷7=4, 4" "80÷4=20, 0"
Replace-Text "95÷9=10, 5" "79÷3=26, 1"
Replace-Text "85÷3=28, 1" "84÷6=14, 0"
Replace-Text "57÷8=7, 1" "91÷6=15, 1"
Replace-Text "35÷5=7, 0" "76÷7=10, 6"
